$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) "The stories can become more descriptive " + "using a persona"
#    -> merge into a single run. Find/Replace across the two runs
#    naturally collapses them into one run in the saved OOXML.
# -------------------------------------------------------------------
$d.Content.Find.Execute(
    "The stories can become more descriptive using a persona", $true,
    $false, $false, $false, $false, $true, 1, $false,
    "The stories can become more descriptive using a persona", 2) | Out-Null

# -------------------------------------------------------------------
# 2) Insert the new "As a.I want. so that." block (with the relocated
#    _GoBack bookmark) plus the four new "User Stories" lines between
#    the "Product owners agreement..." paragraph and the
#    "Given.when.then" heading; then add a second run "." to the
#    "Given.when.then" heading.
# -------------------------------------------------------------------

# Locate the "Product owners agreement..." paragraph.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Product owners agreement*") {
        $anchor = $p
        break
    }
}

# Drop the old _GoBack bookmark that currently sits at the end of
# that paragraph; it gets recreated at the end of the new heading.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# New Heading1: "As a.I want. so that. <<<"
$anchor.Range.InsertParagraphAfter()
$heading = $anchor.Next()
$heading.Range.Text = "As a…I want… so that… <<<"
$heading.Style = "Heading 1"

# Re-add the _GoBack bookmark collapsed right after the new run
# (before the paragraph mark). A bookmark anchored on a truly
# collapsed range at a paragraph's last text position round-trips
# incorrectly, so it is created on a one-character placeholder range
# and the placeholder is removed afterwards, which leaves the
# bookmark collapsed in the right place.
$headingNow = $anchor.Next()
$pos = $headingNow.Range.End - 1
$d.Range($pos, $pos).InsertAfter("X")
$markRange = $d.Range($pos, $pos + 1)
$markRange.Bookmarks.Add("_GoBack") | Out-Null
$d.Range($pos, $pos + 1).Text = ""

# Four NoSpacing lines.
$lines = @(
    "User Stories",
    "As a “someone”",
    "I want to “to do something”",
    "So that “I can achieve a personal goal”"
)
$prev = $anchor.Next()
foreach ($line in $lines) {
    $prev.Range.InsertParagraphAfter()
    $p = $prev.Next()
    $p.Range.Text = $line
    $p.Style = "No Spacing"
    $prev = $p
}

# Append a distinct second run "." to the "Given.when.then" heading.
$given = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Given*when*then*") {
        $given = $p
        break
    }
}
$gPos = $given.Range.End - 1
$gRange = $d.Range($gPos, $gPos)
$gRange.InsertAfter("…")
$newRunRange = $d.Range($gPos, $gPos + 1)
$newRunRange.Bold = 1
$newRunRange.Bold = 0

# -------------------------------------------------------------------
# 3) "Features (Behaviour Driven Development" + ")(" +
#    "Development View)" (split around proofErr tags) -> merge into
#    one run and drop the proofErr markers.
# -------------------------------------------------------------------
$d.Content.Find.Execute(
    "Features (Behaviour Driven Development)(Development View)", $true,
    $false, $false, $false, $false, $true, 1, $false,
    "Features (Behaviour Driven Development)(Development View)", 2) | Out-Null
